# Battery_Data sheet (sheet1): update existing B2:B5 values
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Battery_Data")

$ws1.Range("B2").Value = 26.0753354553
$ws1.Range("B3").Value = 14.080681145862
$ws1.Range("B4").Value = 0.28161362291724
$ws1.Range("B5").Value = 4.97494874334

# Yearly BRC sheet (sheet2): update existing B2:B5 values and append rows 6:21
$ws2 = $wb.Worksheets.Item("Yearly BRC")

$ws2.Range("B2").Value = 0.6286452380429166
$ws2.Range("B3").Value = 0.6433874021842154
$ws2.Range("B4").Value = 0.6606993624998114
$ws2.Range("B5").Value = 0.6806390408598679

# Copy the style of the last existing labeled row (A5) down to the new rows
# so new A6:A21 label cells pick up the same formatting (border/bold/center).
$ws2.Range("A5").Copy()
$ws2.Range("A6:A21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws2.Range("A6").Value = "Battery Replacement Cost at y = 5"
$ws2.Range("B6").Value = 0.7038081714019271

$ws2.Range("A7").Value = "Battery Replacement Cost at y = 6"
$ws2.Range("B7").Value = 0.730752324459822

$ws2.Range("A8").Value = "Battery Replacement Cost at y = 7"
$ws2.Range("B8").Value = 0.7616039398545558

$ws2.Range("A9").Value = "Battery Replacement Cost at y = 8"
$ws2.Range("B9").Value = 0.7965992632711851

$ws2.Range("A10").Value = "Battery Replacement Cost at y = 9"
$ws2.Range("B10").Value = 0.836035758905422

$ws2.Range("A11").Value = "Battery Replacement Cost at y = 10"
$ws2.Range("B11").Value = 0.8803156244852249

$ws2.Range("A12").Value = "Battery Replacement Cost at y = 11"
$ws2.Range("B12").Value = 0.9299165750643577

$ws2.Range("A13").Value = "Battery Replacement Cost at y = 12"
$ws2.Range("B13").Value = 0.9853259660240716

$ws2.Range("A14").Value = "Battery Replacement Cost at y = 13"
$ws2.Range("B14").Value = 1.047170289693914

$ws2.Range("A15").Value = "Battery Replacement Cost at y = 14"
$ws2.Range("B15").Value = 1.11597647473226

$ws2.Range("A16").Value = "Battery Replacement Cost at y = 15"
$ws2.Range("B16").Value = 1.191990881979759

$ws2.Range("A17").Value = "Battery Replacement Cost at y = 16"
$ws2.Range("B17").Value = 1.27511091357169

$ws2.Range("A18").Value = "Battery Replacement Cost at y = 17"
$ws2.Range("B18").Value = 1.364759030235794

$ws2.Range("A19").Value = "Battery Replacement Cost at y = 18"
$ws2.Range("B19").Value = 1.459724038914086

$ws2.Range("A20").Value = "Battery Replacement Cost at y = 19"
$ws2.Range("B20").Value = 1.558547236289639

$ws2.Range("A21").Value = "Battery Replacement Cost at y = 20"
$ws2.Range("B21").Value = 1.65975115358533
